$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.429.00"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "3.316.96"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'558.53"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").Value = "'143.33"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.317.66"
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'27.09"
$ws.Range("E15").Value = "  -4.39%  "
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("D18").Value = "60.379.68"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").Value = "'14.54"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").Value = "'374.95"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "'74.12"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "3.436.39"
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("E27").Value = "  -7.39%  "
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'7.28"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").Value = "'166.79"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "'26.83"
$ws.Range("E41").Value = "  -14.52%  "
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "2.365.44"
$ws.Range("E48").Value = "  -6.89%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("D51").Value = "'21.53"
$ws.Range("E51").Value = "  -4.33%  "
